$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20:116 down to 21:117
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the latest weekly price entry
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44565
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108002
$ws.Range("J20").Value = "Mango"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 450
$ws.Range("N20").Value = 6500
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 6750
$ws.Range("Q20").Value = "$/bandeja 4 kilos"
$ws.Range("R20").Value = "Perú"
$ws.Range("S20").Value = 1688
$ws.Range("T20").Value = 4
